{"js": "// Gameplay_Testing_Document.docx edit:\n// The \"User can exit the game...\" row's video-test reference changes\n// from \"(See Video Test 1)\" to \"(See Video Test 3)\".\n//\n// Locate the paragraph by an unambiguous anchor phrase, then restrict\n// the text search to that paragraph only (there is another, unrelated\n// \"(See Video Test 1)\" earlier in the document that must stay intact).\n\nconst body = context.document.body;\nconst anchor = body.search(\"Goes straight to menu.\", { matchCase: true });\nanchor.load(\"items\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const targetParagraph = anchor.items[0].paragraphs.getFirst();\n\n  const target = targetParagraph.search(\"Video Test 1\", { matchCase: true });\n  target.load(\"items\");\n  await context.sync();\n\n  if (target.items.length > 0) {\n    target.items[0].insertText(\"Video Test 3\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Gameplay_Testing_Document.docx edit:\n# The \"User can exit the game...\" row's video-test reference changes\n# from \"(See Video Test 1)\" to \"(See Video Test 3)\".\n#\n# Locate the paragraph by an unambiguous anchor phrase, then restrict\n# the replacement to that paragraph only (there is another, unrelated\n# \"(See Video Test 1)\" earlier in the document that must stay intact).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Goes straight to menu*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $paraRange = $target.Range\n    $paraText = $paraRange.Text\n    $relIndex = $paraText.IndexOf(\"Video Test 1\")\n\n    if ($relIndex -ge 0) {\n        $absStart = $paraRange.Start + $relIndex\n        $absEnd = $absStart + \"Video Test 1\".Length\n\n        $hit = $d.Range($absStart, $absEnd)\n        $hit.Text = \"Video Test 3\"\n    }\n}\n"}
